$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell I1 ---
# Set the label, then copy the formatting (bold font, border, centered
# alignment) used by the rest of the header row (A1:H1) so the new header
# cell picks up the same style index instead of creating a new one.
$ws.Range("I1").Value = "norm_sample"
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# --- Data cells I2:I96 ---
# Simulated standard-normal sample values ("norm_sample" column) added
# alongside the existing ecdf column.
$normSample = @(
    -0.06884756893673077,
    -0.05747885415475073,
    -0.05409011469717837,
    -0.04893597275464123,
    -0.04270437276595071,
    -0.04264752826659474,
    -0.04045763297614108,
    -0.03435591152980129,
    -0.03290726718297413,
    -0.03104092188805441,
    -0.028556055836795,
    -0.02679669085194228,
    -0.02643806212275866,
    -0.02552513068504525,
    -0.02316557483509203,
    -0.02313972646396837,
    -0.0212821417883994,
    -0.02017259528916622,
    -0.01974285654372802,
    -0.0187100333554594,
    -0.01794010598874173,
    -0.01733577913495614,
    -0.01642362578693973,
    -0.0159266714797445,
    -0.01485750084075787,
    -0.01472186685021784,
    -0.01398161163124536,
    -0.01358951627189297,
    -0.01331279720686511,
    -0.01295535249349514,
    -0.01116226779421779,
    -0.01111284894134258,
    -0.01022384370548639,
    -0.01015481704480267,
    -0.01008915410924837,
    -0.009257830926931814,
    -0.008239367445834081,
    -0.008234710748315281,
    -0.006187766520348362,
    -0.004802629570806549,
    -0.004324588592705947,
    -0.003274872137717096,
    -0.002664264109184237,
    -0.002644876547946708,
    -0.00257317950560736,
    -0.001938412156323156,
    -0.0009821958358117336,
    -0.0004444190928998785,
    0.0003149312196317896,
    0.001734303913876858,
    0.001905627237624808,
    0.002060168879335924,
    0.002785159436380215,
    0.003175512489880537,
    0.003453511046511433,
    0.004233166909151772,
    0.00437837429349911,
    0.006393850650235813,
    0.007030265213488995,
    0.007347836654934577,
    0.007643478670453655,
    0.007726608419276571,
    0.008084802244504889,
    0.008521509237455907,
    0.009073449296860367,
    0.009992845817947507,
    0.01026845215075704,
    0.01099613608185137,
    0.0141101164532052,
    0.01520595310196653,
    0.01545001307148072,
    0.01579339777863938,
    0.01605276762101358,
    0.01631449483921896,
    0.01883576485206128,
    0.01944224292139983,
    0.01947556215657805,
    0.02083281027531672,
    0.02125990775286905,
    0.02204494578394293,
    0.02435477551692228,
    0.0270170088350674,
    0.02842940002001671,
    0.02846059321963981,
    0.02905969938468173,
    0.03014550406742116,
    0.03038788186466854,
    0.03053181655111144,
    0.03062073717900662,
    0.03212137827295825,
    0.03230174570592321,
    0.03241406151680494,
    0.0443311401490801,
    0.04517453033017867,
    0.04715823800072949
)

$startRow = 2
for ($i = 0; $i -lt $normSample.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 9).Value = $normSample[$i]
}
